$wb = $excel.ActiveWorkbook

$wsTemps   = $wb.Worksheets.Item("Temps")
$wsCamions = $wb.Worksheets.Item("Camions")

# --- Sheet "Temps": move the little E1:E8 lookup table one column to the right (E -> F) ---
[void]$wsTemps.Range("E1:E8").Cut($wsTemps.Range("F1:F8"))

# --- Sheet "Camions": add two new columns of data (Cout d'usage / Cout fixe) ---
$wsCamions.Range("D1").Value = "Cout d'usage"
$wsCamions.Range("E1").Value = "Cout fixe"

$wsCamions.Range("D2").Value = 5
$wsCamions.Range("E2").Value = 200

$wsCamions.Range("D3").Value = 4
$wsCamions.Range("E3").Value = 150

$wsCamions.Range("D4").Value = 3
$wsCamions.Range("E4").Value = 100

# --- View state: Camions keeps a selection but is no longer the active tab ---
$wsCamions.Activate()
[void]$wsCamions.Range("E4").Select()

# --- Temps becomes the active tab with the new F2:F8 range selected ---
$wsTemps.Activate()
[void]$wsTemps.Range("F2:F8").Select()
